$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "MOCK_DATA (6)" to "Sheet 1"
$ws.Name = "Sheet 1"

# Widen columns B, C and G to fit their contents (best-fit column widths)
$ws.Columns("B").ColumnWidth = 9.611979166666666
$ws.Columns("C").ColumnWidth = 10.265625
$ws.Columns("G").ColumnWidth = 19.611979166666668

# Update the view: zoom to 200% and move the selection to E14
$excel.ActiveWindow.Zoom = 200
$ws.Range("E14").Select() | Out-Null
